$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Error")

$ws.Range("B28").Value = "JobNotInTown"
$ws.Range("C28").Value = "在城镇中才能切换职业"

$ws.Range("C28").Select() | Out-Null

$tcs = $wb.Theme.ThemeColorScheme
$lt1 = $tcs.Colors(2)
$lt1.RGB = 16777215
